# Form Page and screenshot - 09/29/2021
# Adds a second data row (row 3) to the "in" worksheet, mirroring row 2's
# John Wick sample record except for the Id column which becomes "ID02",
# and re-creates the mailto hyperlink on the new EmailID cell (B3).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of sample data (row 3) - same contents as row 2, new Id.
$ws.Range("A3").Value = "ID02"
$ws.Range("B3").Value = "john@wick.com"
$ws.Range("C3").Value = "John"
$ws.Range("D3").Value = "Wick"
$ws.Range("E3").Value = "John@123"
$ws.Range("F3").Value = "Microsoft"
$ws.Range("G3").Value = "4-41-88/2 sanathnagar"
$ws.Range("H3").Value = "Jagathgirigutta"
$ws.Range("I3").Value = "Hyderabad"
$ws.Range("J3").Value = 50037
$ws.Range("K3").Value = 8019271171

# Recreate the mailto hyperlink on the EmailID cell, matching B2's hyperlink.
$null = $ws.Hyperlinks.Add($ws.Range("B3"), "mailto:john@wick.com")

# Hyperlinks.Add stamps its own style onto the cell; reapply B2's style so
# B3 keeps the same formatting as the existing hyperlink cell.
$ws.Range("B3").Style = $ws.Range("B2").Style

# Final selection left on the sheet after the edit.
$null = $ws.Range("G13").Select()
